$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Long Cypher query text for the new "SamplesTab" row (row 3, column B)
$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE s.study_acronym IN ["D"] 
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@

# Long Cypher query text for the new "FilesTab" row (row 4, column B)
$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE s.study_acronym IN ["D"] 
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@

# Write new rows column-by-column (matches shared-string insertion order of the source edit)
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("C4").Value = $ws.Range("C2").Value2

$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("D4").Value = $ws.Range("D2").Value2

$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("E4").Value = $ws.Range("E2").Value2

# Apply the same "wrap text" style used by B2/C2 to the new long-text cells
$ws.Range("B3:C4").WrapText = $true

# Row heights to roughly match autosized wrapped-text rows
$ws.Rows.Item(2).RowHeight = 316.8
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6

# Column widths (approximate autofit results for the new content; inputs chosen
# so that Excel's internal pixel-snapping lands as close as possible to the
# target stored widths)
$ws.Columns.Item(1).ColumnWidth = 11.91667
$ws.Columns.Item(2).ColumnWidth = 75.25001
$ws.Columns.Item(3).ColumnWidth = 46.91667
$ws.Columns.Item(4).ColumnWidth = 40.91667
$ws.Columns.Item(5).ColumnWidth = 39.58334

# Update the sheet view: zoom level and active selection cell
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("B2").Select() | Out-Null
